# "new wallet info json format"
# The source JSON feeding this sheet now emits two additional wallet
# snapshot dates. Append them to the bottom of the existing "Date" column
# (column A) on Sheet1, as plain text (matching every other date already
# in that column, which are stored as text, not real date values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDates = @("2024-11-15", "2024-08-20")

$lastRow = 38
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $lastRow + 1 + $i
    $cell = $ws.Cells.Item($row, 1)

    # Force text interpretation so Excel doesn't silently coerce the
    # "yyyy-mm-dd" string into a date serial number, then drop the
    # temporary number-format override so the cell ends up with the same
    # (default) style as its neighbours.
    $cell.NumberFormat = "@"
    $cell.Value = $newDates[$i]
    $cell.ClearFormats()
}
